$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("New terms")
$table = $ws.ListObjects.Item("Table3")

# Add a new row to the "Table3" table (rows 2-27 -> now 2-28) and fill it in
$newRow = $table.ListRows.Add()
$ws.Range("A28").Formula = '=IF(LEN(TRIM(B28))=0,0,LEN(TRIM(B28))-LEN(SUBSTITUTE(B28," ",""))+1)'
$ws.Range("B28").Value = "2-axis CNC machine (akin to our robot)"
$ws.Range("C28").Value = "https://youtu.be/05W4egqLVEM"

$null = $ws.Range("B28").Select()
